$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows before the existing row 1053 (2019-11-29), shifting all
# subsequent rows down by 9. This makes room for the missing trading days
# 2019-11-18 .. 2019-11-28 that were added to the price history.
$ws.Rows("1053:1061").Insert()

$newRows = @(
    @{ Row=1053; A=1574035200; B="2019-11-18"; E=1.96; F=1.97; G=1.96; H=1.96; I=576500 },
    @{ Row=1054; A=1574121600; B="2019-11-19"; E=1.95; F=1.98; G=1.94; H=1.98; I=43300 },
    @{ Row=1055; A=1574208000; B="2019-11-20"; E=1.99; F=1.99; G=1.95; H=1.97; I=21300 },
    @{ Row=1056; A=1574294400; B="2019-11-21"; E=1.97; F=1.97; G=1.93; H=1.93; I=44300 },
    @{ Row=1057; A=1574380800; B="2019-11-22"; E=1.96; F=1.96; G=1.91; H=1.95; I=16500 },
    @{ Row=1058; A=1574640000; B="2019-11-25"; E=1.91; F=1.95; G=1.91; H=1.95; I=57800 },
    @{ Row=1059; A=1574726400; B="2019-11-26"; E=1.91; F=1.95; G=1.91; H=1.95; I=21500 },
    @{ Row=1060; A=1574812800; B="2019-11-27"; E=1.95; F=1.95; G=1.91; H=1.92; I=91500 },
    @{ Row=1061; A=1574899200; B="2019-11-28"; E=1.92; F=1.92; G=1.89; H=1.89; I=199100 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Range("A$row").Value = $r.A

    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("B$row").Value = $r.B
    $ws.Range("B$row").Style = "Normal"

    $ws.Range("C$row").NumberFormat = "@"
    $ws.Range("C$row").Value = "5263"
    $ws.Range("C$row").Style = "Normal"

    $ws.Range("D$row").NumberFormat = "@"
    $ws.Range("D$row").Value = "SUNCON"
    $ws.Range("D$row").Style = "Normal"

    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
}
